$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new SD_ext row (row 32) with the extra/external SD values ---
$ws.Range("A32").Value = "SD_ext"
$ws.Range("E32").Value = 0.078665152695345922
$ws.Range("I32").Value = 0.12417619053534083

# --- Update the E and I column formulas (rows 2-26) to fold in the new
#     external SD term E$32 / I$32, combined in quadrature with the
#     original half-width-based estimate ---
for ($r = 2; $r -le 26; $r++) {
    $ws.Range("E$r").Formula = "=SQRT(((D$r-C$r)/2/TINV(0.05,`$X$r)*SQRT(`$X$r))^2+E`$32^2)"
    $ws.Range("I$r").Formula = "=SQRT(((H$r-G$r)/2/TINV(0.05,`$X$r)*SQRT(`$X$r))^2+I`$32^2)"
}

# --- Column widths for E and I (best-fit sized after the new content) ---
$ws.Columns.Item(5).ColumnWidth = 11.1666666666667
$ws.Columns.Item(9).ColumnWidth = 10.1666666666667

# --- Selection moved to Q31 (matches the saved cursor position) ---
$ws.Range("Q31").Select()
